$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.9220969785113576, 0.5763078753274383, 1.421419631572951, 3.465216388719809)"
$ws.Range("C2").Value = "JSU(-1.2363144543294735, 1.3849304519231738, 1.8947561095990895, 6.588223694572884)"
$ws.Range("D2").Value = "NIG(1.028932030541391, 0.7633964303387992, 1.5163135433156785, 2.7553621122380583)"
$ws.Range("E2").Value = "JSU(-1.257614330154591, 1.2963921600481483, 2.323139274040787, 5.093724143633565)"
